$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.498.77"
$ws.Range("E2").Value = "  +2.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.470.15"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9559"
$ws.Range("E5").Value = "  -4.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "280.97"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3706"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3176"
$ws.Range("E8").Value = "  +2.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.80"
$ws.Range("E9").Value = "  +5.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.058"
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06666"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.595"
$ws.Range("E13").Value = "  +3.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.21"
$ws.Range("E14").Value = "  +6.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.239"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.476.16"
$ws.Range("E16").Value = "  +3.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001035"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05724"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9550"
$ws.Range("E19").Value = "  -4.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.08"
$ws.Range("E20").Value = "  -3.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.661"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.69"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.19"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.259"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.723.50"
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.288"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "137.91"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.55"
$ws.Range("E28").Value = "  +3.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.639.81"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.58"
$ws.Range("E30").Value = "  +4.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.957"
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.312"
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8300"
$ws.Range("E33").Value = "  -7.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.621"
$ws.Range("E34").Value = "  +28.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07822"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06033"
$ws.Range("E36").Value = "  +6.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.907"
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.63"
$ws.Range("E38").Value = "  -6.30%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02067"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9704"
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1888"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.330"
$ws.Range("E43").Value = "  -13.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5394"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.587"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.38"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.98"
$ws.Range("E47").Value = "  +11.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5298"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("E50").Value = "  +4.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.049"
$ws.Range("E51").Value = "  -0.16%  "
